$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "57.428.65"
$ws.Range("E2").Value = "  +1.66%  "
$ws.Range("D3").Value = "2.357.34"
$ws.Range("E3").Value = "  +1.06%  "
$ws.Range("E4").Value = "  +0.16%  "
$ws.Range("D5").Value = "'541.33"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +5.39%  "
$ws.Range("D6").Value = "'134.68"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.86%  "
$ws.Range("E7").Value = "  -0.39%  "
$ws.Range("D8").Value = "'0.538"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.80%  "
$ws.Range("D9").Value = "2.357.18"
$ws.Range("E9").Value = "  +0.88%  "
$ws.Range("D10").Value = "'0.102"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.78%  "
$ws.Range("E11").Value = "  +0.93%  "
$ws.Range("D12").Value = "'5.39"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.73%  "
$ws.Range("D13").Value = "'0.354"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +4.52%  "
$ws.Range("D14").Value = "2.749.07"
$ws.Range("E14").Value = "  +0.11%  "
$ws.Range("D15").Value = "'23.47"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.40%  "
$ws.Range("D16").Value = "57.701.63"
$ws.Range("E16").Value = "  +2.17%  "
$ws.Range("E17").Value = "  +0.90%  "
$ws.Range("D18").Value = "2.350.81"
$ws.Range("E18").Value = "  +0.64%  "
$ws.Range("D19").Value = "'10.56"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.08%  "
$ws.Range("D20").Value = "'334.65"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.97%  "
$ws.Range("E21").Value = "  +1.60%  "
$ws.Range("D22").Value = "'6.74"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.61%  "
$ws.Range("D23").Value = "'0.998"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.03%  "
$ws.Range("D24").Value = "'61.82"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.14%  "
$ws.Range("E25").Value = "  +3.95%  "
$ws.Range("D26").Value = "'0.995"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.35%  "
$ws.Range("D27").Value = "'8.41"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.38%  "
$ws.Range("E28").Value = "  +8.70%  "
$ws.Range("E29").Value = "  +4.39%  "
$ws.Range("D30").Value = "'170.25"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.32%  "
$ws.Range("D31").Value = "0.0₃0736"
$ws.Range("E31").Value = "  +2.10%  "
$ws.Range("D32").Value = "'6.16"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.69%  "
$ws.Range("D33").Value = "'18.55"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.95%  "
$ws.Range("E34").Value = "  +15.17%  "
$ws.Range("E35").Value = "  -0.10%  "
$ws.Range("D36").Value = "'0.994"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.36%  "
$ws.Range("E37").Value = "  -0.33%  "
$ws.Range("D38").Value = "'4.14"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +5.39%  "
$ws.Range("E39").Value = "  +3.15%  "
$ws.Range("D40").Value = "'39.34"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.26%  "
$ws.Range("D41").Value = "'150.47"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.52%  "
$ws.Range("D42").Value = "'0.380"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.78%  "
$ws.Range("E43").Value = "  +1.18%  "
$ws.Range("D44").Value = "'285.51"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.10%  "
$ws.Range("E45").Value = "  +6.20%  "
$ws.Range("D47").Value = "'0.0505"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.96%  "
$ws.Range("D48").Value = "'0.562"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.70%  "
$ws.Range("E49").Value = "  +1.96%  "
$ws.Range("E50").Value = "  +2.13%  "
$ws.Range("D51").Value = "'0.382"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.04%  "
